$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 15.33750333333333
$ws.Cells.Item(2, 8).Value = 46.01251
$ws.Cells.Item(2, 9).Value = 0.1440483515229198
$ws.Cells.Item(2, 10).Value = 0.1440483515229198
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 7.731686
$ws.Cells.Item(2, 14).Value = 23.195058
$ws.Cells.Item(2, 15).Value = 0.1963057092861306
$ws.Cells.Item(2, 16).Value = 0.1963057092861306
$ws.Cells.Item(2, 17).Value = 118.5847597972867
$ws.Cells.Item(2, 18).Value = 1067.26283817558
$ws.Cells.Item(2, 19).Value = 0.02827751381720464
$ws.Cells.Item(2, 20).Value = 0.02827751381720464

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 15.33750333333333
$ws.Cells.Item(3, 8).Value = 46.01251
$ws.Cells.Item(3, 9).Value = 0.1440483515229198
$ws.Cells.Item(3, 10).Value = 0.1440483515229198
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 2.796453333333333
$ws.Cells.Item(3, 14).Value = 8.38936
$ws.Cells.Item(3, 15).Value = 0.07100129972758389
$ws.Cells.Item(3, 16).Value = 0.07100129972758387
$ws.Cells.Item(3, 17).Value = 42.89061232151111
$ws.Cells.Item(3, 18).Value = 386.0155108936
$ws.Cells.Item(3, 19).Value = 0.0102276201817432
$ws.Cells.Item(3, 20).Value = 0.01022762018174319

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 15.33750333333333
$ws.Cells.Item(4, 8).Value = 46.01251
$ws.Cells.Item(4, 9).Value = 0.1440483515229198
$ws.Cells.Item(4, 10).Value = 0.1440483515229198
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 25.627053
$ws.Cells.Item(4, 14).Value = 76.881159
$ws.Cells.Item(4, 15).Value = 0.6506649152692259
$ws.Cells.Item(4, 16).Value = 0.6506649152692259
$ws.Cells.Item(4, 17).Value = 393.05501081101
$ws.Cells.Item(4, 18).Value = 3537.49509729909
$ws.Cells.Item(4, 19).Value = 0.0937272084383323
$ws.Cells.Item(4, 20).Value = 0.0937272084383323

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 15.33750333333333
$ws.Cells.Item(5, 8).Value = 46.01251
$ws.Cells.Item(5, 9).Value = 0.1440483515229198
$ws.Cells.Item(5, 10).Value = 0.1440483515229198
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.230753333333334
$ws.Cells.Item(5, 14).Value = 9.692260000000001
$ws.Cells.Item(5, 15).Value = 0.08202807571705974
$ws.Cells.Item(5, 16).Value = 0.08202807571705972
$ws.Cells.Item(5, 17).Value = 49.55169001917778
$ws.Cells.Item(5, 18).Value = 445.9652101726001
$ws.Cells.Item(5, 19).Value = 0.01181600908563971
$ws.Cells.Item(5, 20).Value = 0.0118160090856397

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 64.92210766666668
$ws.Cells.Item(6, 8).Value = 194.766323
$ws.Cells.Item(6, 9).Value = 0.609742171429695
$ws.Cells.Item(6, 10).Value = 0.6097421714296949
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 7.731686
$ws.Cells.Item(6, 14).Value = 23.195058
$ws.Cells.Item(6, 15).Value = 0.1963057092861306
$ws.Cells.Item(6, 16).Value = 0.1963057092861306
$ws.Cells.Item(6, 17).Value = 501.9573509368594
$ws.Cells.Item(6, 18).Value = 4517.616158431734
$ws.Cells.Item(6, 19).Value = 0.1196958694441717
$ws.Cells.Item(6, 20).Value = 0.1196958694441717

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 64.92210766666668
$ws.Cells.Item(7, 8).Value = 194.766323
$ws.Cells.Item(7, 9).Value = 0.609742171429695
$ws.Cells.Item(7, 10).Value = 0.6097421714296949
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.796453333333333
$ws.Cells.Item(7, 14).Value = 8.38936
$ws.Cells.Item(7, 15).Value = 0.07100129972758389
$ws.Cells.Item(7, 16).Value = 0.07100129972758387
$ws.Cells.Item(7, 17).Value = 181.5516443914756
$ws.Cells.Item(7, 18).Value = 1633.96479952328
$ws.Cells.Item(7, 19).Value = 0.04329248667022761
$ws.Cells.Item(7, 20).Value = 0.04329248667022759

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 64.92210766666668
$ws.Cells.Item(8, 8).Value = 194.766323
$ws.Cells.Item(8, 9).Value = 0.609742171429695
$ws.Cells.Item(8, 10).Value = 0.6097421714296949
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 25.627053
$ws.Cells.Item(8, 14).Value = 76.881159
$ws.Cells.Item(8, 15).Value = 0.6506649152692259
$ws.Cells.Item(8, 16).Value = 0.6506649152692259
$ws.Cells.Item(8, 17).Value = 1663.762294045373
$ws.Cells.Item(8, 18).Value = 14973.86064640836
$ws.Cells.Item(8, 19).Value = 0.3967378383093763
$ws.Cells.Item(8, 20).Value = 0.3967378383093762

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 64.92210766666668
$ws.Cells.Item(9, 8).Value = 194.766323
$ws.Cells.Item(9, 9).Value = 0.609742171429695
$ws.Cells.Item(9, 10).Value = 0.6097421714296949
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 3.230753333333334
$ws.Cells.Item(9, 14).Value = 9.692260000000001
$ws.Cells.Item(9, 15).Value = 0.08202807571705974
$ws.Cells.Item(9, 16).Value = 0.08202807571705972
$ws.Cells.Item(9, 17).Value = 209.7473157511089
$ws.Cells.Item(9, 18).Value = 1887.725841759981
$ws.Cells.Item(9, 19).Value = 0.05001597700591944
$ws.Cells.Item(9, 20).Value = 0.05001597700591942

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 10.67805633333333
$ws.Cells.Item(10, 8).Value = 32.034169
$ws.Cells.Item(10, 9).Value = 0.1002872748488753
$ws.Cells.Item(10, 10).Value = 0.1002872748488753
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 7.731686
$ws.Cells.Item(10, 14).Value = 23.195058
$ws.Cells.Item(10, 15).Value = 0.1963057092861306
$ws.Cells.Item(10, 16).Value = 0.1963057092861306
$ws.Cells.Item(10, 17).Value = 82.55937865964465
$ws.Cells.Item(10, 18).Value = 743.0344079368019
$ws.Cells.Item(10, 19).Value = 0.01968696462158158
$ws.Cells.Item(10, 20).Value = 0.01968696462158158

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 10.67805633333333
$ws.Cells.Item(11, 8).Value = 32.034169
$ws.Cells.Item(11, 9).Value = 0.1002872748488753
$ws.Cells.Item(11, 10).Value = 0.1002872748488753
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 2.796453333333333
$ws.Cells.Item(11, 14).Value = 8.38936
$ws.Cells.Item(11, 15).Value = 0.07100129972758389
$ws.Cells.Item(11, 16).Value = 0.07100129972758387
$ws.Cells.Item(11, 17).Value = 29.86068622687111
$ws.Cells.Item(11, 18).Value = 268.74617604184
$ws.Cells.Item(11, 19).Value = 0.007120526860407577
$ws.Cells.Item(11, 20).Value = 0.007120526860407575

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 10.67805633333333
$ws.Cells.Item(12, 8).Value = 32.034169
$ws.Cells.Item(12, 9).Value = 0.1002872748488753
$ws.Cells.Item(12, 10).Value = 0.1002872748488753
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 25.627053
$ws.Cells.Item(12, 14).Value = 76.881159
$ws.Cells.Item(12, 15).Value = 0.6506649152692259
$ws.Cells.Item(12, 16).Value = 0.6506649152692259
$ws.Cells.Item(12, 17).Value = 273.647115591319
$ws.Cells.Item(12, 18).Value = 2462.824040321871
$ws.Cells.Item(12, 19).Value = 0.06525341119212498
$ws.Cells.Item(12, 20).Value = 0.06525341119212498

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 10.67805633333333
$ws.Cells.Item(13, 8).Value = 32.034169
$ws.Cells.Item(13, 9).Value = 0.1002872748488753
$ws.Cells.Item(13, 10).Value = 0.1002872748488753
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 3.230753333333334
$ws.Cells.Item(13, 14).Value = 9.692260000000001
$ws.Cells.Item(13, 15).Value = 0.08202807571705974
$ws.Cells.Item(13, 16).Value = 0.08202807571705972
$ws.Cells.Item(13, 17).Value = 34.49816609243778
$ws.Cells.Item(13, 18).Value = 310.48349483194
$ws.Cells.Item(13, 19).Value = 0.00822637217476112
$ws.Cells.Item(13, 20).Value = 0.00822637217476112

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 15.537021
$ws.Cells.Item(14, 8).Value = 46.611063
$ws.Cells.Item(14, 9).Value = 0.14592220219851
$ws.Cells.Item(14, 10).Value = 0.14592220219851
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 7.731686
$ws.Cells.Item(14, 14).Value = 23.195058
$ws.Cells.Item(14, 15).Value = 0.1963057092861306
$ws.Cells.Item(14, 16).Value = 0.1963057092861306
$ws.Cells.Item(14, 17).Value = 120.127367747406
$ws.Cells.Item(14, 18).Value = 1081.146309726654
$ws.Cells.Item(14, 19).Value = 0.02864536140317267
$ws.Cells.Item(14, 20).Value = 0.02864536140317266

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 15.537021
$ws.Cells.Item(15, 8).Value = 46.611063
$ws.Cells.Item(15, 9).Value = 0.14592220219851
$ws.Cells.Item(15, 10).Value = 0.14592220219851
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 2.796453333333333
$ws.Cells.Item(15, 14).Value = 8.38936
$ws.Cells.Item(15, 15).Value = 0.07100129972758389
$ws.Cells.Item(15, 16).Value = 0.07100129972758387
$ws.Cells.Item(15, 17).Value = 43.44855416552001
$ws.Cells.Item(15, 18).Value = 391.03698748968
$ws.Cells.Item(15, 19).Value = 0.01036066601520551
$ws.Cells.Item(15, 20).Value = 0.01036066601520551

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 15.537021
$ws.Cells.Item(16, 8).Value = 46.611063
$ws.Cells.Item(16, 9).Value = 0.14592220219851
$ws.Cells.Item(16, 10).Value = 0.14592220219851
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 25.627053
$ws.Cells.Item(16, 14).Value = 76.881159
$ws.Cells.Item(16, 15).Value = 0.6506649152692259
$ws.Cells.Item(16, 16).Value = 0.6506649152692259
$ws.Cells.Item(16, 17).Value = 398.168060629113
$ws.Cells.Item(16, 18).Value = 3583.512545662017
$ws.Cells.Item(16, 19).Value = 0.09494645732939236
$ws.Cells.Item(16, 20).Value = 0.09494645732939234

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 15.537021
$ws.Cells.Item(17, 8).Value = 46.611063
$ws.Cells.Item(17, 9).Value = 0.14592220219851
$ws.Cells.Item(17, 10).Value = 0.14592220219851
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 3.230753333333334
$ws.Cells.Item(17, 14).Value = 9.692260000000001
$ws.Cells.Item(17, 15).Value = 0.08202807571705974
$ws.Cells.Item(17, 16).Value = 0.08202807571705972
$ws.Cells.Item(17, 17).Value = 50.19628238582001
$ws.Cells.Item(17, 18).Value = 451.7665414723801
$ws.Cells.Item(17, 19).Value = 0.01196971745073948
$ws.Cells.Item(17, 20).Value = 0.01196971745073948
